$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the cryptos list refresh.
# Force text number format so values like "1.000" or "27.313.76" are not
# auto-converted to numbers by Excel, matching the original inline-string cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.313.76"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.903.27"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.36"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5221"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3779"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07297"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.32"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08237"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +7.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "97.06"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.906.07"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008650"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.337.40"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.103"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.70"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.446"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.311"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.53"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.41%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.746"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.50"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.854"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.939"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09249"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8012"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.239"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.947"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.596"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5730"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02010"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.03%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.038"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.584"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.07"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.96%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4896"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.13"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.631"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.06"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.86"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05948"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.48%  "
